# Hortaliza, Vega Modelo de Temuco - Achicoria: add a new weekly price
# record. This inserts a brand-new row at row 60 (pushing the existing
# rows 60-99 down to 61-100, which also grows the sheet's used range
# from A1:R99 to A1:R100) and fills it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(60).Insert()

$ws.Cells.Item(60, 1).Value = 10
$ws.Cells.Item(60, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(60, 3).Value = 'La Araucanía'
$ws.Cells.Item(60, 4).Value = 45062
$ws.Cells.Item(60, 5).Value = 9
$ws.Cells.Item(60, 6).Value = 100112010
$ws.Cells.Item(60, 7).Value = 'Achicoria'
$ws.Cells.Item(60, 8).Value = 'Sin especificar'
$ws.Cells.Item(60, 9).Value = 'Primera'
$ws.Cells.Item(60, 10).Value = 110
$ws.Cells.Item(60, 11).Value = 10000
$ws.Cells.Item(60, 12).Value = 10000
$ws.Cells.Item(60, 13).Value = 10000
$ws.Cells.Item(60, 14).Value = '$/caja 18 unidades'
$ws.Cells.Item(60, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(60, 16).Value = 556
$ws.Cells.Item(60, 17).Value = 18
$ws.Cells.Item(60, 18).Value = 'Hortaliza'
